$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID Resources-HCP")

# Clear rows 100-108 entirely (values, hyperlinks, formatting, row height)
$ws.Range("A100:G108").EntireRow.Clear()
$ws.Range("A100:G108").EntireRow.AutoFit()

# Restore the "hyperlink-look" formatting (style used by G97:G99) on the now-empty G column cells
$ws.Range("G99").Copy()
$ws.Range("G100:G108").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the frozen-pane view state
$ws.Application.ActiveWindow.ScrollRow = 95
$ws.Range("C103").Select()
